$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 372.1085406666666
$ws.Range("H2").Value = 1116.325622
$ws.Range("I2").Value = 0.8095247142929753
$ws.Range("J2").Value = 0.8095247142929753
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 163.7119853333333
$ws.Range("N2").Value = 491.135956
$ws.Range("O2").Value = 0.2754003062401033
$ws.Range("P2").Value = 0.2754003062401033
$ws.Range("Q2").Value = 60918.6279520294
$ws.Range("R2").Value = 548267.6515682645
$ws.Range("S2").Value = 0.2229433542252175
$ws.Range("T2").Value = 0.2229433542252175

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 372.1085406666666
$ws.Range("H3").Value = 1116.325622
$ws.Range("I3").Value = 0.8095247142929753
$ws.Range("J3").Value = 0.8095247142929753
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 172.558497
$ws.Range("N3").Value = 517.675491
$ws.Range("O3").Value = 0.290282124557779
$ws.Range("P3").Value = 0.290282124557779
$ws.Range("Q3").Value = 64210.49049830336
$ws.Range("R3").Value = 577894.4144847302
$ws.Range("S3").Value = 0.2349905539469939
$ws.Range("T3").Value = 0.2349905539469939

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 372.1085406666666
$ws.Range("H4").Value = 1116.325622
$ws.Range("I4").Value = 0.8095247142929753
$ws.Range("J4").Value = 0.8095247142929753
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.4881643333333
$ws.Range("N4").Value = 334.464493
$ws.Range("O4").Value = 0.1875481171218523
$ws.Range("P4").Value = 0.1875481171218523
$ws.Range("Q4").Value = 41485.69813168218
$ws.Range("R4").Value = 373371.2831851396
$ws.Range("S4").Value = 0.1518248359292529
$ws.Range("T4").Value = 0.151824835929253

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 372.1085406666666
$ws.Range("H5").Value = 1116.325622
$ws.Range("I5").Value = 0.8095247142929753
$ws.Range("J5").Value = 0.8095247142929753
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 86.95798233333333
$ws.Range("N5").Value = 260.873947
$ws.Range("O5").Value = 0.1462828449356383
$ws.Range("P5").Value = 0.1462828449356383
$ws.Range("Q5").Value = 32357.80790537444
$ws.Range("R5").Value = 291220.27114837
$ws.Range("S5").Value = 0.1184195782524862
$ws.Range("T5").Value = 0.1184195782524862

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 372.1085406666666
$ws.Range("H6").Value = 1116.325622
$ws.Range("I6").Value = 0.8095247142929753
$ws.Range("J6").Value = 0.8095247142929753
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 59.73436333333333
$ws.Range("N6").Value = 179.20309
$ws.Range("O6").Value = 0.100486607144627
$ws.Range("P6").Value = 0.100486607144627
$ws.Range("Q6").Value = 22227.6667676191
$ws.Range("R6").Value = 200049.0009085719
$ws.Range("S6").Value = 0.0813463919390246
$ws.Range("T6").Value = 0.0813463919390246

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 69.70494733333334
$ws.Range("H7").Value = 209.114842
$ws.Range("I7").Value = 0.1516435969830949
$ws.Range("J7").Value = 0.1516435969830949
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.7119853333333
$ws.Range("N7").Value = 491.135956
$ws.Range("O7").Value = 0.2754003062401033
$ws.Range("P7").Value = 0.2754003062401033
$ws.Range("Q7").Value = 11411.53531549544
$ws.Range("R7").Value = 102703.8178394589
$ws.Range("S7").Value = 0.04176269304849514
$ws.Range("T7").Value = 0.04176269304849514

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 69.70494733333334
$ws.Range("H8").Value = 209.114842
$ws.Range("I8").Value = 0.1516435969830949
$ws.Range("J8").Value = 0.1516435969830949
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 172.558497
$ws.Range("N8").Value = 517.675491
$ws.Range("O8").Value = 0.290282124557779
$ws.Range("P8").Value = 0.290282124557779
$ws.Range("Q8").Value = 12028.18094530416
$ws.Range("R8").Value = 108253.6285077374
$ws.Range("S8").Value = 0.04401942550783638
$ws.Range("T8").Value = 0.04401942550783638

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 69.70494733333334
$ws.Range("H9").Value = 209.114842
$ws.Range("I9").Value = 0.1516435969830949
$ws.Range("J9").Value = 0.1516435969830949
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.4881643333333
$ws.Range("N9").Value = 334.464493
$ws.Range("O9").Value = 0.1875481171218523
$ws.Range("P9").Value = 0.1875481171218523
$ws.Range("Q9").Value = 7771.276623145012
$ws.Range("R9").Value = 69941.48960830511
$ws.Range("S9").Value = 0.02844047108776444
$ws.Range("T9").Value = 0.02844047108776445

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 69.70494733333334
$ws.Range("H10").Value = 209.114842
$ws.Range("I10").Value = 0.1516435969830949
$ws.Range("J10").Value = 0.1516435969830949
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 86.95798233333333
$ws.Range("N10").Value = 260.873947
$ws.Range("O10").Value = 0.1462828449356383
$ws.Range("P10").Value = 0.1462828449356383
$ws.Range("Q10").Value = 6061.401578757931
$ws.Range("R10").Value = 54552.61420882137
$ws.Range("S10").Value = 0.0221828567829605
$ws.Range("T10").Value = 0.0221828567829605

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 69.70494733333334
$ws.Range("H11").Value = 209.114842
$ws.Range("I11").Value = 0.1516435969830949
$ws.Range("J11").Value = 0.1516435969830949
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.73436333333333
$ws.Range("N11").Value = 179.20309
$ws.Range("O11").Value = 0.100486607144627
$ws.Range("P11").Value = 0.100486607144627
$ws.Range("Q11").Value = 4163.780650140197
$ws.Range("R11").Value = 37474.02585126177
$ws.Range("S11").Value = 0.01523815055603839
$ws.Range("T11").Value = 0.01523815055603839

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.387356
$ws.Range("H12").Value = 1.162068
$ws.Range("I12").Value = 0.0008426956679571845
$ws.Range("J12").Value = 0.0008426956679571845
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 163.7119853333333
$ws.Range("N12").Value = 491.135956
$ws.Range("O12").Value = 0.2754003062401033
$ws.Range("P12").Value = 0.2754003062401033
$ws.Range("Q12").Value = 63.41481979077867
$ws.Range("R12").Value = 570.733378117008
$ws.Range("S12").Value = 0.000232078645022617
$ws.Range("T12").Value = 0.000232078645022617

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.387356
$ws.Range("H13").Value = 1.162068
$ws.Range("I13").Value = 0.0008426956679571845
$ws.Range("J13").Value = 0.0008426956679571845
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 172.558497
$ws.Range("N13").Value = 517.675491
$ws.Range("O13").Value = 0.290282124557779
$ws.Range("P13").Value = 0.290282124557779
$ws.Range("Q13").Value = 66.84156916393201
$ws.Range("R13").Value = 601.574122475388
$ws.Range("S13").Value = 0.0002446194888502482
$ws.Range("T13").Value = 0.0002446194888502482

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.387356
$ws.Range("H14").Value = 1.162068
$ws.Range("I14").Value = 0.0008426956679571845
$ws.Range("J14").Value = 0.0008426956679571845
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 111.4881643333333
$ws.Range("N14").Value = 334.464493
$ws.Range("O14").Value = 0.1875481171218523
$ws.Range("P14").Value = 0.1875481171218523
$ws.Range("Q14").Value = 43.18560938350267
$ws.Range("R14").Value = 388.670484451524
$ws.Range("S14").Value = 0.0001580459858321116
$ws.Range("T14").Value = 0.0001580459858321116

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.387356
$ws.Range("H15").Value = 1.162068
$ws.Range("I15").Value = 0.0008426956679571845
$ws.Range("J15").Value = 0.0008426956679571845
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 86.95798233333333
$ws.Range("N15").Value = 260.873947
$ws.Range("O15").Value = 0.1462828449356383
$ws.Range("P15").Value = 0.1462828449356383
$ws.Range("Q15").Value = 33.68369620471067
$ws.Range("R15").Value = 303.153265842396
$ws.Range("S15").Value = 0.000123271919723715
$ws.Range("T15").Value = 0.000123271919723715

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.387356
$ws.Range("H16").Value = 1.162068
$ws.Range("I16").Value = 0.0008426956679571845
$ws.Range("J16").Value = 0.0008426956679571845
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 59.73436333333333
$ws.Range("N16").Value = 179.20309
$ws.Range("O16").Value = 0.100486607144627
$ws.Range("P16").Value = 0.100486607144627
$ws.Range("Q16").Value = 23.13846404334667
$ws.Range("R16").Value = 208.24617639012
$ws.Range("S16").Value = 0.00008467962852849262
$ws.Range("T16").Value = 0.00008467962852849262

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.425022
$ws.Range("H17").Value = 1.275066
$ws.Range("I17").Value = 0.0009246383125251667
$ws.Range("J17").Value = 0.0009246383125251667
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 163.7119853333333
$ws.Range("N17").Value = 491.135956
$ws.Range("O17").Value = 0.2754003062401033
$ws.Range("P17").Value = 0.2754003062401033
$ws.Range("Q17").Value = 69.581195430344
$ws.Range("R17").Value = 626.230758873096
$ws.Range("S17").Value = 0.0002546456744307633
$ws.Range("T17").Value = 0.0002546456744307633

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.425022
$ws.Range("H18").Value = 1.275066
$ws.Range("I18").Value = 0.0009246383125251667
$ws.Range("J18").Value = 0.0009246383125251667
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 172.558497
$ws.Range("N18").Value = 517.675491
$ws.Range("O18").Value = 0.290282124557779
$ws.Range("P18").Value = 0.290282124557779
$ws.Range("Q18").Value = 73.341157511934
$ws.Range("R18").Value = 660.070417607406
$ws.Range("S18").Value = 0.000268405973807325
$ws.Range("T18").Value = 0.000268405973807325

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.425022
$ws.Range("H19").Value = 1.275066
$ws.Range("I19").Value = 0.0009246383125251667
$ws.Range("J19").Value = 0.0009246383125251667
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 111.4881643333333
$ws.Range("N19").Value = 334.464493
$ws.Range("O19").Value = 0.1875481171218523
$ws.Range("P19").Value = 0.1875481171218523
$ws.Range("Q19").Value = 47.384922581282
$ws.Range("R19").Value = 426.464303231538
$ws.Range("S19").Value = 0.0001734141745328218
$ws.Range("T19").Value = 0.0001734141745328219

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 0.425022
$ws.Range("H20").Value = 1.275066
$ws.Range("I20").Value = 0.0009246383125251667
$ws.Range("J20").Value = 0.0009246383125251667
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 86.95798233333333
$ws.Range("N20").Value = 260.873947
$ws.Range("O20").Value = 0.1462828449356383
$ws.Range("P20").Value = 0.1462828449356383
$ws.Range("Q20").Value = 36.959055567278
$ws.Range("R20").Value = 332.631500105502
$ws.Range("S20").Value = 0.0001352587228926693
$ws.Range("T20").Value = 0.0001352587228926693

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 0.425022
$ws.Range("H21").Value = 1.275066
$ws.Range("I21").Value = 0.0009246383125251667
$ws.Range("J21").Value = 0.0009246383125251667
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 59.73436333333333
$ws.Range("N21").Value = 179.20309
$ws.Range("O21").Value = 0.100486607144627
$ws.Range("P21").Value = 0.100486607144627
$ws.Range("Q21").Value = 25.38841857266
$ws.Range("R21").Value = 228.49576715394
$ws.Range("S21").Value = 0.00009291376686158725
$ws.Range("T21").Value = 0.00009291376686158725

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 17.03711166666666
$ws.Range("H22").Value = 51.111335
$ws.Range("I22").Value = 0.03706435474344739
$ws.Range("J22").Value = 0.03706435474344739
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 163.7119853333333
$ws.Range("N22").Value = 491.135956
$ws.Range("O22").Value = 0.2754003062401033
$ws.Range("P22").Value = 0.2754003062401033
$ws.Range("Q22").Value = 2789.179375295695
$ws.Range("R22").Value = 25102.61437766126
$ws.Range("S22").Value = 0.01020753464693724
$ws.Range("T22").Value = 0.01020753464693724

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 17.03711166666666
$ws.Range("H23").Value = 51.111335
$ws.Range("I23").Value = 0.03706435474344739
$ws.Range("J23").Value = 0.03706435474344739
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 172.558497
$ws.Range("N23").Value = 517.675491
$ws.Range("O23").Value = 0.290282124557779
$ws.Range("P23").Value = 0.290282124557779
$ws.Range("Q23").Value = 2939.898382421165
$ws.Range("R23").Value = 26459.08544179048
$ws.Range("S23").Value = 0.0107591196402911
$ws.Range("T23").Value = 0.0107591196402911

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 17.03711166666666
$ws.Range("H24").Value = 51.111335
$ws.Range("I24").Value = 0.03706435474344739
$ws.Range("J24").Value = 0.03706435474344739
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 111.4881643333333
$ws.Range("N24").Value = 334.464493
$ws.Range("O24").Value = 0.1875481171218523
$ws.Range("P24").Value = 0.1875481171218523
$ws.Range("Q24").Value = 1899.436305258684
$ws.Range("R24").Value = 17094.92674732815
$ws.Range("S24").Value = 0.006951349944469953
$ws.Range("T24").Value = 0.006951349944469954

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 17.03711166666666
$ws.Range("H25").Value = 51.111335
$ws.Range("I25").Value = 0.03706435474344739
$ws.Range("J25").Value = 0.03706435474344739
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 86.95798233333333
$ws.Range("N25").Value = 260.873947
$ws.Range("O25").Value = 0.1462828449356383
$ws.Range("P25").Value = 0.1462828449356383
$ws.Range("Q25").Value = 1481.512855321027
$ws.Range("R25").Value = 13333.61569788924
$ws.Range("S25").Value = 0.005421879257575205
$ws.Range("T25").Value = 0.005421879257575205

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 17.03711166666666
$ws.Range("H26").Value = 51.111335
$ws.Range("I26").Value = 0.03706435474344739
$ws.Range("J26").Value = 0.03706435474344739
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 59.73436333333333
$ws.Range("N26").Value = 179.20309
$ws.Range("O26").Value = 0.100486607144627
$ws.Range("P26").Value = 0.100486607144627
$ws.Range("Q26").Value = 1017.701018447239
$ws.Range("R26").Value = 9159.309166025148
$ws.Range("S26").Value = 0.003724471254173889
$ws.Range("T26").Value = 0.003724471254173889
